$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 52
$ws.Range("H52").Value = 849.6667
$ws.Range("I52").Value = 999.5
$ws.Range("J52").Value = 550
$ws.Range("K52").Value = 2998.5
$ws.Range("L52").Value = 1650
$ws.Range("M52").Value = -2838.5
$ws.Range("N52").Value = -1970

# Row 125
$ws.Range("H125").Value = 2340.5
$ws.Range("I125").Value = 1932
$ws.Range("J125").Value = 2476.6667
$ws.Range("K125").Value = 17388
$ws.Range("L125").Value = 22290.0003
$ws.Range("M125").Value = -14928
$ws.Range("N125").Value = -27210.0003

# Row 137
$ws.Range("H137").Value = 2942
$ws.Range("I137").Value = 1949.1305
$ws.Range("J137").Value = 5018
$ws.Range("K137").Value = 5847.3915
$ws.Range("L137").Value = 15054
$ws.Range("M137").Value = -3297.3915
$ws.Range("N137").Value = -20154

# Row 138
$ws.Range("H138").Value = 3628.8474
$ws.Range("I138").Value = 1874.5385
$ws.Range("J138").Value = 4124.6304
$ws.Range("K138").Value = 5623.6155
$ws.Range("L138").Value = 12373.8912
$ws.Range("M138").Value = -483.6154999999999
$ws.Range("N138").Value = -22653.8912

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 17600.555
$ws.Range("I2").Value = 6276
$ws.Range("J2").Value = 31756.25
$ws.Range("K2").Value = 6276
$ws.Range("L2").Value = 31756.25
$ws.Range("M2").Value = -6163

# Row 32
$ws.Range("H32").Value = 1239.8948
$ws.Range("I32").Value = 1215.6857
$ws.Range("J32").Value = 1522.3334
$ws.Range("K32").Value = 1215.6857
$ws.Range("L32").Value = 1522.3334
$ws.Range("M32").Value = -928.6857

# Row 61
$ws.Range("H61").Value = 5877.2593
$ws.Range("I61").Value = 5213.0454
$ws.Range("J61").Value = 8799.799999999999
$ws.Range("K61").Value = 5213.0454
$ws.Range("L61").Value = 8799.799999999999
$ws.Range("M61").Value = -5001.0454

# Row 116
$ws.Range("H116").Value = 17600.555
$ws.Range("I116").Value = 6276
$ws.Range("J116").Value = 31756.25
$ws.Range("K116").Value = 6276
$ws.Range("L116").Value = 31756.25
$ws.Range("M116").Value = -3982

# Row 132
$ws.Range("H132").Value = 4246.7427
$ws.Range("I132").Value = 3311.1475
$ws.Range("J132").Value = 10588
$ws.Range("K132").Value = 9933.442500000001
$ws.Range("L132").Value = 31764
$ws.Range("M132").Value = -7403.442500000001

# Row 136
$ws.Range("H136").Value = 5877.2593
$ws.Range("I136").Value = 5213.0454
$ws.Range("J136").Value = 8799.799999999999
$ws.Range("K136").Value = 15639.1362
$ws.Range("L136").Value = 26399.4
$ws.Range("M136").Value = -13089.1362

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 17600.555
$ws.Range("I3").Value = 6276
$ws.Range("J3").Value = 31756.25
$ws.Range("K3").Value = 6276
$ws.Range("L3").Value = 31756.25
$ws.Range("M3").Value = -6162

# Row 86
$ws.Range("H86").Value = 4005.96
$ws.Range("I86").Value = 2333.842
$ws.Range("J86").Value = 9301
$ws.Range("K86").Value = 2333.842
$ws.Range("L86").Value = 9301
$ws.Range("M86").Value = -1210.842

# Row 89
$ws.Range("H89").Value = 4005.96
$ws.Range("I89").Value = 2333.842
$ws.Range("J89").Value = 9301
$ws.Range("K89").Value = 11669.21
$ws.Range("L89").Value = 46505
$ws.Range("M89").Value = -6053.210000000001

# Row 134
$ws.Range("H134").Value = 2615.1777
$ws.Range("I134").Value = 1928.8611
$ws.Range("J134").Value = 5360.4443
$ws.Range("K134").Value = 5786.5833
$ws.Range("L134").Value = 16081.3329
$ws.Range("M134").Value = -3251.5833

# Row 140
$ws.Range("H140").Value = 65999.2
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 65999.2
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 65999.2
$ws.Range("N140").Value = -76359.2

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 24935.72
$ws.Range("I31").Value = 2958.85
$ws.Range("J31").Value = 39586.965
$ws.Range("K31").Value = 2958.85
$ws.Range("L31").Value = 39586.965
$ws.Range("M31").Value = -2663.85
$ws.Range("N31").Value = -40176.965

# Row 34
$ws.Range("H34").Value = 24935.72
$ws.Range("I34").Value = 2958.85
$ws.Range("J34").Value = 39586.965
$ws.Range("K34").Value = 2958.85
$ws.Range("L34").Value = 39586.965
$ws.Range("M34").Value = -2756.85
$ws.Range("N34").Value = -39990.965

# Row 39
$ws.Range("H39").Value = 9389.166999999999
$ws.Range("I39").Value = 4886.5
$ws.Range("J39").Value = 18394.5
$ws.Range("K39").Value = 4886.5
$ws.Range("L39").Value = 18394.5
$ws.Range("M39").Value = -4495.5

# Row 49
$ws.Range("H49").Value = 9389.166999999999
$ws.Range("I49").Value = 4886.5
$ws.Range("J49").Value = 18394.5
$ws.Range("K49").Value = 4886.5
$ws.Range("L49").Value = 18394.5
$ws.Range("M49").Value = -4704.5

# Row 58
$ws.Range("H58").Value = 5367.933
$ws.Range("I58").Value = 3086.9333
$ws.Range("J58").Value = 7648.933
$ws.Range("K58").Value = 3086.9333
$ws.Range("L58").Value = 7648.933
$ws.Range("M58").Value = -2883.9333

# Row 62
$ws.Range("H62").Value = 8062
$ws.Range("I62").Value = 5586.4
$ws.Range("J62").Value = 14251
$ws.Range("K62").Value = 5586.4
$ws.Range("L62").Value = 14251
$ws.Range("M62").Value = -4962.4

# Row 65
$ws.Range("H65").Value = 8062
$ws.Range("I65").Value = 5586.4
$ws.Range("J65").Value = 14251
$ws.Range("K65").Value = 27932
$ws.Range("L65").Value = 71255
$ws.Range("M65").Value = -24812

# Row 99
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()

# Row 107
$ws.Range("H107").Value = 2734.818
$ws.Range("I107").Value = 1218
$ws.Range("J107").Value = 3998.8333
$ws.Range("K107").Value = 1218
$ws.Range("L107").Value = 3998.8333
$ws.Range("M107").Value = 702
$ws.Range("N107").Value = -7838.8333

# Row 122
$ws.Range("H122").Value = 5208.375
$ws.Range("I122").Value = 2774.65
$ws.Range("J122").Value = 9264.583000000001
$ws.Range("K122").Value = 8323.950000000001
$ws.Range("L122").Value = 27793.749
$ws.Range("M122").Value = -5873.950000000001
$ws.Range("N122").Value = -32693.749

# Row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

# Row 136
$ws.Range("H136").Value = 5367.933
$ws.Range("I136").Value = 3086.9333
$ws.Range("J136").Value = 7648.933
$ws.Range("K136").Value = 9260.7999
$ws.Range("L136").Value = 22946.799
$ws.Range("M136").Value = -6710.7999

# Row 140
$ws.Range("H140").Value = 133999.2
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 133999.2
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 133999.2
$ws.Range("N140").Value = -144359.2

$ws = $wb.Worksheets.Item("CUL")
# Row 38
$ws.Range("H38").Value = 3921.6365
$ws.Range("I38").Value = 102.4
$ws.Range("J38").Value = 7104.3335
$ws.Range("K38").Value = 307.2
$ws.Range("L38").Value = 21313.0005
$ws.Range("M38").Value = 39.79999999999995

# Row 68
$ws.Range("H68").Value = 5275.5317
$ws.Range("I68").Value = 999.75
$ws.Range("J68").Value = 5673.2793
$ws.Range("K68").Value = 2999.25
$ws.Range("L68").Value = 17019.8379
$ws.Range("M68").Value = -2188.25
$ws.Range("N68").Value = -18641.8379

# Row 71
$ws.Range("H71").Value = 5275.5317
$ws.Range("I71").Value = 999.75
$ws.Range("J71").Value = 5673.2793
$ws.Range("K71").Value = 8997.75
$ws.Range("L71").Value = 51059.5137
$ws.Range("M71").Value = -4941.75
$ws.Range("N71").Value = -59171.5137

# Row 86
$ws.Range("H86").Value = 1050
$ws.Range("I86").Value = 416.66666
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 1249.99998
$ws.Range("L86").Value = 6000
$ws.Range("M86").Value = -63.99998000000005
$ws.Range("N86").Value = -8372

# Row 89
$ws.Range("H89").Value = 1050
$ws.Range("I89").Value = 416.66666
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 3749.99994
$ws.Range("L89").Value = 18000
$ws.Range("M89").Value = 2178.00006
$ws.Range("N89").Value = -29856

$ws = $wb.Worksheets.Item("GSM")
# Row 14
$ws.Range("H14").Value = 854.1
$ws.Range("I14").Value = 598.875
$ws.Range("J14").Value = 1875
$ws.Range("K14").Value = 598.875
$ws.Range("L14").Value = 1875
$ws.Range("M14").Value = -430.875
$ws.Range("N14").Value = -2211

# Row 33
$ws.Range("H33").Value = 53333
$ws.Range("I33").Value = 49999
$ws.Range("J33").Value = 55000
$ws.Range("K33").Value = 49999
$ws.Range("L33").Value = 55000
$ws.Range("M33").Value = -49747
$ws.Range("N33").Value = -55504

# Row 46
$ws.Range("H46").Value = 1190.625
$ws.Range("I46").Value = 1289.2858
$ws.Range("J46").Value = 500
$ws.Range("K46").Value = 1289.2858
$ws.Range("L46").Value = 500
$ws.Range("M46").Value = -1133.2858

# Row 104
$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

# Row 133
$ws.Range("H133").Value = 69996.664
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 69996.664
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 69996.664
$ws.Range("N133").Value = -80116.664

$ws = $wb.Worksheets.Item("LTW")
# Row 110
$ws.Range("H110").Value = 41644
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 41644
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 41644
$ws.Range("N110").Value = -49824

# Row 122
$ws.Range("H122").Value = 7495.1113
$ws.Range("I122").Value = 5909.3335
$ws.Range("J122").Value = 10666.667
$ws.Range("K122").Value = 17728.0005
$ws.Range("L122").Value = 32000.001
$ws.Range("M122").Value = -15278.0005

# Row 132
$ws.Range("H132").Value = 5703.6665
$ws.Range("I132").Value = 3993.818
$ws.Range("J132").Value = 8390.571
$ws.Range("K132").Value = 11981.454
$ws.Range("L132").Value = 25171.713
$ws.Range("M132").Value = -9451.454000000002
$ws.Range("N132").Value = -30231.713

$ws = $wb.Worksheets.Item("WVR")
# Row 41
$ws.Range("H41").Value = 19040.334
$ws.Range("I41").Value = 17342
$ws.Range("J41").Value = 19889.5
$ws.Range("K41").Value = 17342
$ws.Range("L41").Value = 19889.5
$ws.Range("M41").Value = -16952
$ws.Range("N41").Value = -20669.5

# Row 59
$ws.Range("H59").Value = 43500
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 43500
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 43500
$ws.Range("N59").Value = -44976

# Row 81
$ws.Range("H81").Value = 15768
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 15768
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 31536
$ws.Range("N81").Value = -33658

# Row 84
$ws.Range("H84").Value = 15768
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 15768
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 157680
$ws.Range("N84").Value = -168288

# Row 126
$ws.Range("H126").Value = 15002.667
$ws.Range("I126").Value = 10004
$ws.Range("J126").Value = 25000
$ws.Range("K126").Value = 30012
$ws.Range("L126").Value = 75000
$ws.Range("M126").Value = -27542
$ws.Range("N126").Value = -79940
